{"js": "const replacements = [\n  [\"2025-05-21 Wednesday\", \"2025-05-22 Thursday\"],\n  [\"614\u00f75=122, 4\", \"523\u00f78=65, 3\"],\n  [\"127\u00f72=63, 1\", \"571\u00f73=190, 1\"],\n  [\"681\u00f74=170, 1\", \"262\u00f75=52, 2\"],\n  [\"130\u00f76=21, 4\", \"845\u00f76=140, 5\"],\n  [\"773\u00f73=257, 2\", \"362\u00f78=45, 2\"],\n  [\"699\u00f75=139, 4\", \"733\u00f72=366, 1\"],\n  [\"669\u00f74=167, 1\", \"991\u00f79=110, 1\"],\n  [\"818\u00f77=116, 6\", \"988\u00f79=109, 7\"],\n  [\"899\u00f78=112, 3\", \"746\u00f74=186, 2\"],\n  [\"288\u00f79=32, 0\", \"431\u00f72=215, 1\"],\n  [\"951\u00f79=105, 6\", \"200\u00f79=22, 2\"],\n  [\"774\u00f73=258, 0\", \"172\u00f75=34, 2\"],\n  [\"461\u00f79=51, 2\", \"762\u00f75=152, 2\"],\n  [\"478\u00f74=119, 2\", \"521\u00f74=130, 1\"],\n  [\"574\u00f78=71, 6\", \"528\u00f78=66, 0\"],\n  [\"935\u00f79=103, 8\", \"996\u00f73=332, 0\"],\n  [\"181\u00f76=30, 1\", \"395\u00f73=131, 2\"],\n  [\"849\u00f79=94, 3\", \"189\u00f76=31, 3\"],\n  [\"169\u00f78=21, 1\", \"128\u00f79=14, 2\"],\n  [\"523\u00f76=87, 1\", \"539\u00f72=269, 1\"],\n  [\"986\u00f73=328, 2\", \"270\u00f74=67, 2\"],\n  [\"707\u00f78=88, 3\", \"109\u00f78=13, 5\"],\n  [\"761\u00f75=152, 1\", \"475\u00f79=52, 7\"],\n  [\"238\u00f73=79, 1\", \"287\u00f75=57, 2\"],\n  [\"224\u00f72=112, 0\", \"669\u00f75=133, 4\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-05-21 Wednesday\", \"2025-05-22 Thursday\"),\n  @(\"614\u00f75=122, 4\", \"523\u00f78=65, 3\"),\n  @(\"127\u00f72=63, 1\", \"571\u00f73=190, 1\"),\n  @(\"681\u00f74=170, 1\", \"262\u00f75=52, 2\"),\n  @(\"130\u00f76=21, 4\", \"845\u00f76=140, 5\"),\n  @(\"773\u00f73=257, 2\", \"362\u00f78=45, 2\"),\n  @(\"699\u00f75=139, 4\", \"733\u00f72=366, 1\"),\n  @(\"669\u00f74=167, 1\", \"991\u00f79=110, 1\"),\n  @(\"818\u00f77=116, 6\", \"988\u00f79=109, 7\"),\n  @(\"899\u00f78=112, 3\", \"746\u00f74=186, 2\"),\n  @(\"288\u00f79=32, 0\", \"431\u00f72=215, 1\"),\n  @(\"951\u00f79=105, 6\", \"200\u00f79=22, 2\"),\n  @(\"774\u00f73=258, 0\", \"172\u00f75=34, 2\"),\n  @(\"461\u00f79=51, 2\", \"762\u00f75=152, 2\"),\n  @(\"478\u00f74=119, 2\", \"521\u00f74=130, 1\"),\n  @(\"574\u00f78=71, 6\", \"528\u00f78=66, 0\"),\n  @(\"935\u00f79=103, 8\", \"996\u00f73=332, 0\"),\n  @(\"181\u00f76=30, 1\", \"395\u00f73=131, 2\"),\n  @(\"849\u00f79=94, 3\", \"189\u00f76=31, 3\"),\n  @(\"169\u00f78=21, 1\", \"128\u00f79=14, 2\"),\n  @(\"523\u00f76=87, 1\", \"539\u00f72=269, 1\"),\n  @(\"986\u00f73=328, 2\", \"270\u00f74=67, 2\"),\n  @(\"707\u00f78=88, 3\", \"109\u00f78=13, 5\"),\n  @(\"761\u00f75=152, 1\", \"475\u00f79=52, 7\"),\n  @(\"238\u00f73=79, 1\", \"287\u00f75=57, 2\"),\n  @(\"224\u00f72=112, 0\", \"669\u00f75=133, 4\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Find/Replace failed for: $oldText\"\n  }\n}"}
